$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a day-month-year style label (e.g. "08-10-2021") as plain
# literal text, without letting Excel's input parser reinterpret it as a
# date serial number and without leaving any new cell-style behind.
# We do this by putting the literal text into the cell via a formula
# (="08-10-2021") and then collapsing the formula down to its computed
# string value with Copy + PasteSpecial(values only). The cached string
# is preserved as-is (no date coercion), and PasteSpecial-values does not
# introduce any new number-format / style entries.
function Set-TextLabel($cell, [string]$text) {
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# Row 194: 08-10-2021
Set-TextLabel $ws.Cells.Item(194, 1) "08-10-2021"
$ws.Cells.Item(194, 3).Value = 0.32
$ws.Cells.Item(194, 4).Value = 2.69
$ws.Cells.Item(194, 5).Value = 3.33
$ws.Cells.Item(194, 6).Value = 3.5
$ws.Cells.Item(194, 7).Value = 3.6

# Row 195: 12-10-2021 (no value in column F)
Set-TextLabel $ws.Cells.Item(195, 1) "12-10-2021"
$ws.Cells.Item(195, 3).Value = 0.6
$ws.Cells.Item(195, 4).Value = 2.88
$ws.Cells.Item(195, 5).Value = 3.52
$ws.Cells.Item(195, 7).Value = 3.88

# Row 196: 13-10-2021
Set-TextLabel $ws.Cells.Item(196, 1) "13-10-2021"
$ws.Cells.Item(196, 3).Value = 0.62
$ws.Cells.Item(196, 4).Value = 2.92
$ws.Cells.Item(196, 5).Value = 3.58
$ws.Cells.Item(196, 6).Value = 3.8
$ws.Cells.Item(196, 7).Value = 3.9

# Row 197: 14-10-2021 (only columns A, C, D, E populated)
Set-TextLabel $ws.Cells.Item(197, 1) "14-10-2021"
$ws.Cells.Item(197, 3).Value = 0.54
$ws.Cells.Item(197, 4).Value = 2.75
$ws.Cells.Item(197, 5).Value = 3.15
